# Paar refactorings, en wat commentaar toegevoegd.
# Adds a new worked example (compound-interest) block to the "Testen F#"
# sheet: header row 51 (hoofdsom / rente / looptijd), data row 52
# (principal, rate as %, term via formula "=4+(5/12)" with a literal-text
# comment cell), and the compounded result in row 53.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testen F#")

# --- Row 52: data first, so the shared-string table picks up the
#     "4 + (5/12)" comment text before the row-51 header labels (matches
#     the order the strings were actually typed in by the author). ---
$ws.Range("A52").Value = 5500
$ws.Range("D52").Value = "'4 + (5/12)"
$ws.Range("C52").Formula = "=4+(5/12)"
$ws.Range("B52").Value = 0.038
$ws.Range("B52").NumberFormat = "0.00%"

# --- Row 51: header labels above the data just entered. ---
$ws.Range("A51").Value = "hoofdsom"
$ws.Range("B51").Value = "rente"
$ws.Range("C51").Value = "looptijd"

# --- Row 53: the compounded-value formula. ---
$ws.Range("A53").Formula = "=A52*(1 + B52)^C52"

# --- View state: scroll down and select B52, mirroring the author's
#     on-screen position after finishing the edit. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B52").Select()
